# This edit re-shuffles the species-observation data (columns A,B,D,E,F,G,H,K,Q,R)
# across data rows 5-17 of the "Artfynd" sheet so that each row ends up holding a
# different observation's id/taxon info while location columns (P,S,T,U,V,W,...)
# stay put. K is only populated ("i frukt") for the row that now holds the
# "Skogsknipprot" / Epipactis helleborine observation (row 12); it is cleared
# everywhere else, including row 15 which previously held it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 111943980
$ws.Range("B5").Value = 89317
$ws.Range("D5").Value = 'LC'
$ws.Range("E5").Value = 3215
$ws.Range("F5").Value = 'Rödgul trumpetsvamp'
$ws.Range("G5").Value = 'Craterellus lutescens'
$ws.Range("H5").Value = '(Fr.) Fr.'
$ws.Range("K5").ClearContents()
$ws.Range("Q5").Value = 682877
$ws.Range("R5").Value = 6694410

# Row 6
$ws.Range("A6").Value = 111943992
$ws.Range("B6").Value = 89317
$ws.Range("D6").Value = 'LC'
$ws.Range("E6").Value = 3215
$ws.Range("F6").Value = 'Rödgul trumpetsvamp'
$ws.Range("G6").Value = 'Craterellus lutescens'
$ws.Range("H6").Value = '(Fr.) Fr.'
$ws.Range("K6").ClearContents()
$ws.Range("Q6").Value = 682867
$ws.Range("R6").Value = 6694644

# Row 7
$ws.Range("A7").Value = 111943979
$ws.Range("B7").Value = 96625
$ws.Range("D7").Value = 'LC'
$ws.Range("E7").Value = 504
$ws.Range("F7").Value = 'Guckusko'
$ws.Range("G7").Value = 'Cypripedium calceolus'
$ws.Range("H7").Value = 'L.'
$ws.Range("K7").ClearContents()
$ws.Range("Q7").Value = 682879
$ws.Range("R7").Value = 6694407

# Row 8
$ws.Range("A8").Value = 111943996
$ws.Range("B8").Value = 90466
$ws.Range("D8").Value = 'LC'
$ws.Range("E8").Value = 4769
$ws.Range("F8").Value = 'Svavelriska'
$ws.Range("G8").Value = 'Lactarius scrobiculatus'
$ws.Range("H8").Value = '(Scop.:Fr.) Fr.'
$ws.Range("K8").ClearContents()
$ws.Range("Q8").Value = 682785
$ws.Range("R8").Value = 6694547

# Row 9
$ws.Range("A9").Value = 111943990
$ws.Range("B9").Value = 102166
$ws.Range("D9").Value = 'LC'
$ws.Range("E9").Value = 222412
$ws.Range("F9").Value = 'Tibast'
$ws.Range("G9").Value = 'Daphne mezereum'
$ws.Range("H9").Value = 'L.'
$ws.Range("K9").ClearContents()
$ws.Range("Q9").Value = 682930
$ws.Range("R9").Value = 6694720

# Row 10
$ws.Range("A10").Value = 111943995
$ws.Range("B10").Value = 89033
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 3286
$ws.Range("F10").Value = 'Flattoppad klubbsvamp'
$ws.Range("G10").Value = 'Clavariadelphus truncatus'
$ws.Range("H10").Value = '(Quél.) Donk'
$ws.Range("K10").ClearContents()
$ws.Range("Q10").Value = 682779
$ws.Range("R10").Value = 6694551

# Row 11
$ws.Range("A11").Value = 111943984
$ws.Range("B11").Value = 99850
$ws.Range("D11").Value = 'LC'
$ws.Range("E11").Value = 221235
$ws.Range("F11").Value = 'Vårärt'
$ws.Range("G11").Value = 'Lathyrus vernus'
$ws.Range("H11").Value = '(L.) Bernh.'
$ws.Range("K11").ClearContents()
$ws.Range("Q11").Value = 682929
$ws.Range("R11").Value = 6694685

# Row 12
$ws.Range("A12").Value = 111943997
$ws.Range("B12").Value = 96698
$ws.Range("D12").Value = 'LC'
$ws.Range("E12").Value = 219798
$ws.Range("F12").Value = 'Skogsknipprot'
$ws.Range("G12").Value = 'Epipactis helleborine'
$ws.Range("H12").Value = '(L.) Crantz'
$ws.Range("K12").Value = 'i frukt'
$ws.Range("Q12").Value = 682781
$ws.Range("R12").Value = 6694488

# Row 13
$ws.Range("A13").Value = 111943981
$ws.Range("B13").Value = 96625
$ws.Range("D13").Value = 'LC'
$ws.Range("E13").Value = 504
$ws.Range("F13").Value = 'Guckusko'
$ws.Range("G13").Value = 'Cypripedium calceolus'
$ws.Range("H13").Value = 'L.'
$ws.Range("K13").ClearContents()
$ws.Range("Q13").Value = 682877
$ws.Range("R13").Value = 6694410

# Row 14
$ws.Range("A14").Value = 111943999
$ws.Range("B14").Value = 99850
$ws.Range("D14").Value = 'LC'
$ws.Range("E14").Value = 221235
$ws.Range("F14").Value = 'Vårärt'
$ws.Range("G14").Value = 'Lathyrus vernus'
$ws.Range("H14").Value = '(L.) Bernh.'
$ws.Range("K14").ClearContents()
$ws.Range("Q14").Value = 682757
$ws.Range("R14").Value = 6694406

# Row 15
$ws.Range("A15").Value = 111943983
$ws.Range("B15").Value = 90812
$ws.Range("D15").Value = 'LC'
$ws.Range("E15").Value = 4366
$ws.Range("F15").Value = 'Skarp dropptaggsvamp'
$ws.Range("G15").Value = 'Hydnellum peckii'
$ws.Range("H15").Value = 'Banker'
$ws.Range("K15").ClearContents()
$ws.Range("Q15").Value = 682871
$ws.Range("R15").Value = 6694481

# Row 16
$ws.Range("A16").Value = 111943988
$ws.Range("B16").Value = 107547
$ws.Range("D16").Value = 'NT'
$ws.Range("E16").Value = 220320
$ws.Range("F16").Value = 'Ängsskära'
$ws.Range("G16").Value = 'Serratula tinctoria'
$ws.Range("H16").Value = 'L.'
$ws.Range("K16").ClearContents()
$ws.Range("Q16").Value = 682930
$ws.Range("R16").Value = 6694720

# Row 17
$ws.Range("A17").Value = 111943998
$ws.Range("B17").Value = 98961
$ws.Range("D17").Value = 'LC'
$ws.Range("E17").Value = 222498
$ws.Range("F17").Value = 'Blåsippa'
$ws.Range("G17").Value = 'Hepatica nobilis'
$ws.Range("H17").Value = 'Schreb.'
$ws.Range("K17").ClearContents()
$ws.Range("Q17").Value = 682757
$ws.Range("R17").Value = 6694406
